$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - App. Android-Text
$ws.Range("C2").Value = 1799033.0
$ws.Range("D2").Value = 8177.0
$ws.Range("E2").Value = 114.530322

# Row 3 - App. iOS-Text
$ws.Range("C3").Value = 1269916.0
$ws.Range("D3").Value = 3619.0
$ws.Range("E3").Value = 67.913795
